# Auto-generated Excel COM-interop script to apply golem profit sheet updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (ALC)
$ws.Range("H28").Value = 1056.2727
$ws.Range("I28").Value = 1056.2727
$ws.Range("K28").Value = 1056.2727
$ws.Range("M28").Value = -571.2727

# Row 40 (ALC)
$ws.Range("H40").Value = 1699.6666
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1699.6666
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1699.6666
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -2049.6666

# Row 138 (ALC)
$ws.Range("H138").Value = 3588.7778
$ws.Range("I138").Value = 824.75
$ws.Range("J138").Value = 5800
$ws.Range("K138").Value = 2474.25
$ws.Range("L138").Value = 17400
$ws.Range("M138").Value = 2665.75
$ws.Range("N138").Value = -27680

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 439.06668
$ws.Range("I2").Value = 399
$ws.Range("K2").Value = 399
$ws.Range("M2").Value = -286

# Row 88 (ARM)
$ws.Range("H88").Value = 2068.9333
$ws.Range("I88").Value = 1320
$ws.Range("K88").Value = 1320
$ws.Range("M88").Value = -914

# Row 91 (ARM)
$ws.Range("H91").Value = 2068.9333
$ws.Range("I91").Value = 1320
$ws.Range("K91").Value = 1320
$ws.Range("M91").Value = 84

# Row 116 (ARM)
$ws.Range("H116").Value = 439.06668
$ws.Range("I116").Value = 399
$ws.Range("K116").Value = 399
$ws.Range("M116").Value = 1895

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 439.06668
$ws.Range("I3").Value = 399
$ws.Range("K3").Value = 399
$ws.Range("M3").Value = -285

# Row 105 (BSM)
$ws.Range("H105").Value = 786.5
$ws.Range("I105").Value = 617.5
$ws.Range("J105").Value = 1124.5
$ws.Range("K105").Value = 617.5
$ws.Range("L105").Value = 1124.5
$ws.Range("M105").Value = 1129.5
$ws.Range("N105").Value = -4618.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 13683.25
$ws.Range("I31").Value = 5571.2856
$ws.Range("K31").Value = 5571.2856
$ws.Range("M31").Value = -5276.2856

# Row 34 (CRP)
$ws.Range("H34").Value = 13683.25
$ws.Range("I34").Value = 5571.2856
$ws.Range("K34").Value = 5571.2856
$ws.Range("M34").Value = -5369.2856

# Row 50 (CRP)
$ws.Range("H50").Value = 29450
$ws.Range("I50").Value = 5000
$ws.Range("J50").Value = 53900
$ws.Range("K50").Value = 5000
$ws.Range("L50").Value = 53900
$ws.Range("M50").Value = -4375
$ws.Range("N50").Value = -55150

# Row 59 (CRP)
$ws.Range("H59").Value = 47552
$ws.Range("J59").Value = 65000
$ws.Range("L59").Value = 65000
$ws.Range("N59").Value = -67290

# Row 122 (CRP)
$ws.Range("H122").Value = 1168.25
$ws.Range("I122").Value = 1168.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3504.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1054.75
$ws.Range("N122").Value = $null

# Row 134 (CRP)
$ws.Range("H134").Value = 1418.4286
$ws.Range("I134").Value = 1425.8
$ws.Range("K134").Value = 4277.4
$ws.Range("M134").Value = -1742.4

$ws = $wb.Worksheets.Item("CUL")
# Row 31 (CUL)
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1212

# Row 68 (CUL)
$ws.Range("H68").Value = 2149.2
$ws.Range("J68").Value = 2149.2
$ws.Range("L68").Value = 6447.599999999999
$ws.Range("N68").Value = -8069.599999999999

# Row 71 (CUL)
$ws.Range("H71").Value = 2149.2
$ws.Range("J71").Value = 2149.2
$ws.Range("L71").Value = 19342.8
$ws.Range("N71").Value = -27454.8

# Row 104 (CUL)
$ws.Range("H104").Value = 399.5
$ws.Range("I104").Value = 399.5
$ws.Range("K104").Value = 1198.5
$ws.Range("M104").Value = 1422.5

# Row 122 (CUL)
$ws.Range("H122").Value = 402
$ws.Range("I122").Value = 104
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 936
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = 1514
$ws.Range("N122").Value = -11200

# Row 132 (CUL)
$ws.Range("H132").Value = 250
$ws.Range("I132").Value = 125
$ws.Range("K132").Value = 1125
$ws.Range("M132").Value = 1405

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Range("H122").Value = 6211
$ws.Range("I122").Value = 7000
$ws.Range("J122").Value = 5816.5
$ws.Range("K122").Value = 21000
$ws.Range("L122").Value = 17449.5
$ws.Range("M122").Value = -18550
$ws.Range("N122").Value = -22349.5

# Row 126 (GSM)
$ws.Range("H126").Value = 14999.75
$ws.Range("I126").Value = 9999
$ws.Range("K126").Value = 29997
$ws.Range("M126").Value = -27527

# Row 132 (GSM)
$ws.Range("H132").Value = 1099.3334
$ws.Range("I132").Value = 985
$ws.Range("K132").Value = 2955
$ws.Range("M132").Value = -425

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 500
$ws.Range("K7").Value = 500
$ws.Range("M7").Value = -388

# Row 16 (LTW)
$ws.Range("H16").Value = 1662.625
$ws.Range("I16").Value = 1662.625
$ws.Range("K16").Value = 1662.625
$ws.Range("M16").Value = -1492.625

# Row 43 (LTW)
$ws.Range("H43").Value = 591283.6
$ws.Range("J43").Value = 591283.6
$ws.Range("L43").Value = 591283.6
$ws.Range("N43").Value = -591669.6

# Row 46 (LTW)
$ws.Range("H46").Value = 403518.8
$ws.Range("I46").Value = 1001000
$ws.Range("K46").Value = 1001000
$ws.Range("M46").Value = -1000812

# Row 61 (LTW)
$ws.Range("H61").Value = 998.5
$ws.Range("I61").Value = 998.5
$ws.Range("K61").Value = 998.5
$ws.Range("M61").Value = -796.5

# Row 68 (LTW)
$ws.Range("H68").Value = 3500
$ws.Range("I68").Value = 3500
$ws.Range("K68").Value = 3500
$ws.Range("M68").Value = -2751

# Row 71 (LTW)
$ws.Range("H71").Value = 3500
$ws.Range("I71").Value = 3500
$ws.Range("K71").Value = 17500
$ws.Range("M71").Value = -13756

# Row 113 (LTW)
$ws.Range("H113").Value = 998.5
$ws.Range("I113").Value = 998.5
$ws.Range("K113").Value = 998.5
$ws.Range("M113").Value = 1171.5

# Row 126 (LTW)
$ws.Range("H126").Value = 500
$ws.Range("I126").Value = 500
$ws.Range("K126").Value = 1500
$ws.Range("M126").Value = 970

# Row 136 (LTW)
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 1166.6666
$ws.Range("I81").Value = 1250
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 2500
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -1439
$ws.Range("N81").Value = -4122

# Row 84 (WVR)
$ws.Range("H84").Value = 1166.6666
$ws.Range("I84").Value = 1250
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 12500
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -7196
$ws.Range("N84").Value = -20608

# Row 113 (WVR)
$ws.Range("H113").Value = 4484.154
$ws.Range("I113").Value = 208.27272
$ws.Range("K113").Value = 624.81816
$ws.Range("M113").Value = 1545.18184

# Row 132 (WVR)
$ws.Range("H132").Value = 1199.6666
$ws.Range("I132").Value = 799.5
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2398.5
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = 131.5
$ws.Range("N132").Value = -11060

